$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text (matches source formatting)
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply updated cell values
$ws.Range('D2').Value = '42.711.16'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '2.537.34'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '312.34'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').Value = '100.16'
$ws.Range('D7').Value = '0.565'
$ws.Range('E7').Value = '  -0.79%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('D10').Value = '35.42'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '7.27'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').Value = '2.927.47'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '15.41'
$ws.Range('E15').Value = '  -3.37%  '
$ws.Range('D16').Value = '2.479.33'
$ws.Range('E16').Value = '  -3.60%  '
$ws.Range('D17').Value = '0.813'
$ws.Range('E17').Value = '  -3.05%  '
$ws.Range('D18').Value = '42.710.29'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '6.70'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = '12.29'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = '0.0₃0949'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '69.67'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').Value = '242.73'
$ws.Range('D24').Value = '2.86'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  -2.51%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '25.60'
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('D29').Value = '10.10'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -4.51%  '
$ws.Range('D31').Value = '5.88'
$ws.Range('E31').Value = '  +2.95%  '
$ws.Range('D32').Value = '157.24'
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '2.67'
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.0791'
$ws.Range('E34').Value = '  -1.70%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '3.16'
$ws.Range('E35').Value = '  -3.89%  '
$ws.Range('B36').Value = 'ApeXProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D36').Value = '2.55'
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('D37').Value = '17.76'
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('D38').Value = '1.96'
$ws.Range('E38').Value = '  -5.56%  '
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '4.13'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = '21.73'
$ws.Range('E42').Value = '  -3.07%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = '3.31'
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0298'
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('D46').Value = '1.994.42'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '9.06'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('D48').Value = '2.782.04'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = '0.190'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '79.54'
$ws.Range('E50').Value = '  -2.00%  '
$ws.Range('D51').Value = '72.10'
$ws.Range('E51').Value = '  -1.51%  '
